$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.588.35'
$ws.Range("E2").Value = '  +2.36%  '
$ws.Range("D3").Value = '2.961.31'
$ws.Range("E3").Value = '  +1.13%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.98'
$ws.Range("E5").Value = '  +0.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.15'
$ws.Range("E6").Value = '  -0.15%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '2.960.59'
$ws.Range("E8").Value = '  +1.11%  '
$ws.Range("E9").Value = '  +0.22%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.38'
$ws.Range("E10").Value = '  +5.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.142'
$ws.Range("E11").Value = '  -0.61%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.447'
$ws.Range("E12").Value = '  +1.57%  '
$ws.Range("E13").Value = '  +3.16%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.36'
$ws.Range("E14").Value = '  -1.25%  '
$ws.Range("E15").Value = '  -0.29%  '
$ws.Range("D16").Value = '3.450.13'
$ws.Range("E16").Value = '  +1.07%  '
$ws.Range("D17").Value = '62.512.31'
$ws.Range("E17").Value = '  +2.30%  '
$ws.Range("E18").Value = '  -0.23%  '
$ws.Range("D19").Value = '2.959.23'
$ws.Range("E19").Value = '  +0.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '440.72'
$ws.Range("E20").Value = '  +1.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.42'
$ws.Range("E21").Value = '  -0.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.671'
$ws.Range("E22").Value = '  -1.12%  '
$ws.Range("E23").Value = '  -0.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.84'
$ws.Range("E24").Value = '  +0.43%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.14'
$ws.Range("E25").Value = '  +0.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.96'
$ws.Range("E26").Value = '  +0.62%  '
$ws.Range("E27").Value = '  -3.96%  '
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.60'
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.08'
$ws.Range("E30").Value = '  +1.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.11'
$ws.Range("E31").Value = '  -6.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.56'
$ws.Range("E32").Value = '  -0.47%  '
$ws.Range("E33").Value = '  -2.44%  '
$ws.Range("E34").Value = '  +0.10%  '
$ws.Range("D35").Value = '0.0₃0874'
$ws.Range("E35").Value = '  +0.67%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.992'
$ws.Range("E36").Value = '  -1.77%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.64'
$ws.Range("E37").Value = '  -0.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.05'
$ws.Range("E38").Value = '  +2.68%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '49.61'
$ws.Range("E39").Value = '  -0.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.91'
$ws.Range("E40").Value = '  -3.19%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.55'
$ws.Range("E41").Value = '  -0.44%  '
$ws.Range("E42").Value = '  -5.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.281'
$ws.Range("E43").Value = '  -1.73%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.01'
$ws.Range("E44").Value = '  -7.47%  '
$ws.Range("D45").Value = '2.717.29'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '135.59'
$ws.Range("E46").Value = '  +1.64%  '
$ws.Range("E47").Value = '  -2.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '361.80'
$ws.Range("E48").Value = '  -4.02%  '
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("E50").Value = '  -0.24%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.89'
$ws.Range("E51").Value = '  -4.72%  '
